$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ran E5 CRM test to see variability from probe opened 2021-02-28 (dmb).
# Append the new batch of CRM-accuracy readings below the existing data.
$rows = @(
    @{ Row = 3;  A = 20210228; B = 2211.1392636523301; C = 2234.0700000000002; E = 141; F = "CRM opened 20210216" },
    @{ Row = 4;  A = 20210228; B = 2213.5992626297498; C = 2234.0700000000002; E = 141; F = "CRM opened 20210216" },
    @{ Row = 5;  A = 20210228; B = 2218.8110320871701; C = 2234.0700000000002; E = 141; F = "CRM opened 20210228" },
    @{ Row = 6;  A = 20210228; B = 2208.8304250821302; C = 2234.0700000000002; E = 141; F = "CRM opened 20210228" },
    @{ Row = 7;  A = 20210228; B = 2203.4973775284002; C = 2234.0700000000002; E = 141; F = "CRM opened 20210228" },
    @{ Row = 8;  A = 20210228; B = 2213.3549756320699; C = 2234.0700000000002; E = 141; F = "CRM opened 20210228" },
    @{ Row = 9;  A = 20210228; B = 2230.9452295190799; C = 2234.0700000000002; E = 141; F = "CRM opened 20210228" },
    @{ Row = 10; A = 20210228; B = 2222.5457646464001; C = 2234.0700000000002; E = 141; F = "CRM opened 20210228" }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value = $r.A
    $ws.Cells.Item($i, 2).Value = $r.B
    $ws.Cells.Item($i, 3).Value = $r.C
    $ws.Cells.Item($i, 4).Formula = "=100*(B$i-C$i)/C$i"
    $ws.Cells.Item($i, 5).Value = $r.E
    $ws.Cells.Item($i, 6).Value = $r.F
}

# Matches the final click position recorded in the saved worksheet view.
$ws.Range("H8").Select()
